# Tidsplan_xjobb.xlsx edit
# "continued struggle with making db connection to work"
#
# - Row 15 (2020-01-30 task): plan date pushed a day -> 2020-01-31,
#   hours spent bumped 25 -> 28, and the comment moves from row 14's
#   shared "Plus a buch..." text onto row 15 with the typo fixed
#   ("buch" -> "bunch").
# - All the downstream SUM/ratio formulas (S15, P36, S36, P37, S37,
#   P38, P43, P44, P45) recalc automatically once P15 changes.
# - Active cell/selection moves from Q15 to P15, and the sheet is
#   scrolled over a couple of columns (G1 -> K1 top-left).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P15").Value = 28
$ws.Range("Q15").Value = "2020-01-31 Pending"
$ws.Range("T15").Value = "Plus a bunch of off office hours, problem making connection to db work"

# Move the selection to P15 (was Q15) and scroll the view.
[void]$ws.Range("P15").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 11
$win.ScrollRow = 1
